$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 7) describing the "TbAccount" table, mirroring
# the existing rows (TbGlobalConfig / TbItem / TbShopItem).
$ws.Range("B7").Value = "TbAccount"
$ws.Range("C7").Value = "Account"
# Leading apostrophe forces this to be stored as literal text "False"
# rather than being auto-converted into a Boolean FALSE value.
$ws.Range("D7").Value = "'False"
$ws.Range("E7").Value = "account.xlsx"
$ws.Range("F7").Value = "id"
$ws.Range("G7").Value = "map"
$ws.Range("H7").Value = "c;s"
$ws.Range("I7").Value = "账号测试表"

# Normalize the row's formatting back to the default/normal style so the
# quote-prefix flag picked up from the apostrophe above doesn't linger.
$ws.Range("B7:I7").Style = "Normal"
